# Append 10 new driver records (rows 8-17) to the roster on Sheet1,
# matching the Firstname/Lastname/Role columns already present, then
# left-align the new Role cells (column C) and select C13 to match
# the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDrivers = @(
    @("Shaine",     "Alenin",     "Volunteer"),
    @("Thurston",   "Wayon",      "Volunteer"),
    @("Darnall",    "Frear",      "Permanent"),
    @("Ronalda",    "Carlyon",    "Volunteer"),
    @("Blakeley",   "Gunby",      "Permanent"),
    @("Yorke",      "Hartington", "Volunteer"),
    @("Ryann",      "Britcher",   "Volunteer"),
    @("Adamo",      "Paxton",     "Permanent"),
    @("Kristoffer", "Pagan",      "Permanent"),
    @("Trevor",     "Ollin",      "Volunteer")
)

$row = 8
foreach ($driver in $newDrivers) {
    $ws.Cells.Item($row, 1).Value = $driver[0]
    $ws.Cells.Item($row, 2).Value = $driver[1]
    $ws.Cells.Item($row, 3).Value = $driver[2]

    # Left-align the Role column for the newly added rows.
    $ws.Cells.Item($row, 3).HorizontalAlignment = -4131

    $row++
}

# Column D ends up sized like column C after the new data settles in.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Final selection left on the sheet by the author.
$ws.Range("C13").Select()
